$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-28 Friday", "2025-03-01 Saturday"),
    @("42÷5=", "64÷3="),
    @("42÷9=", "25÷3="),
    @("45÷6=", "62÷3="),
    @("88÷4=", "99÷3="),
    @("69÷8=", "30÷7="),
    @("24÷7=", "84÷6="),
    @("67÷9=", "37÷8="),
    @("21÷8=", "65÷2="),
    @("52÷3=", "37÷7="),
    @("55÷7=", "98÷3="),
    @("45÷5=", "77÷4="),
    @("33÷7=", "70÷3="),
    @("85÷4=", "64÷2="),
    @("89÷6=", "46÷4="),
    @("59÷2=", "44÷5="),
    @("85÷8=", "53÷7="),
    @("38÷5=", "49÷2="),
    @("63÷8=", "46÷7="),
    @("41÷9=", "39÷7="),
    @("44÷2=", "13÷3="),
    @("41÷3=", "10÷8="),
    @("21÷6=", "32÷8="),
    @("69÷3=", "46÷5="),
    @("68÷3=", "28÷9="),
    @("82÷8=", "89÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
